$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("D5").Value = -14.4
$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("C8").ClearContents()
$ws.Range("F8").Value = 17.05
$ws.Range("D9").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("D10").ClearContents()
$ws.Range("F10").Value = 16.43
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = 17.65
$ws.Range("C12").Value = 12.5
$ws.Range("E12").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("E14").Value = -5.4
$ws.Range("C17").Value = 11.2
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("C18").Value = 11.5
$ws.Range("F18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("E21").Value = -8.7
$ws.Range("E22").Value = -6.1
$ws.Range("C23").Value = 12.2
$ws.Range("F23").ClearContents()
$ws.Range("D24").Value = -13.9
$ws.Range("F24").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").ClearContents()
$ws.Range("F26").Value = 17.38
$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("F27").Value = 17
$ws.Range("A28").Value = "SC 105"
$ws.Range("C28").Value = 11.1
$ws.Range("D28").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("A29").Value = "SC 119"
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").ClearContents()
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18
$ws.Range("A32").Value = "SC 193"
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
$ws.Rows("34:35").Delete()

Write-Output "edit applied"
